$wb = $excel.ActiveWorkbook

# --- Sheet "powiat krakowski": flip two "Aktywne" (active) flags to FALSE ---
$wsK = $wb.Worksheets.Item("powiat krakowski")
$wsK.Range("H22").Value = $false
$wsK.Range("H60").Value = $false

# --- Sheet "powiat wielicki": append a new listing as row 55 ---
$wsW = $wb.Worksheets.Item("powiat wielicki")

$wsW.Range("A55").Value = "Atrakcyjna działka budowlana Grajów Wieliczka."
$wsW.Range("B55").Value = "Grajów, Wieliczka, wielicki, małopolskie"
$wsW.Range("C55").Value = 230000

# Dates are stored as plain text (not Excel date serials) in this workbook,
# so force text entry with a leading apostrophe and strip the format Excel
# auto-applies when it recognises the date-like text.
$wsW.Range("D55").Value = "'2025-07-27"
$wsW.Range("D55").Style = "Normal"
$wsW.Range("E55").Value = "'2025-07-27"
$wsW.Range("E55").Style = "Normal"

$wsW.Range("F55").Value = 230000
$wsW.Range("G55").Value = 0.49
$wsW.Range("H55").Value = $true
$wsW.Range("I55").Value = "https://www.otodom.pl/pl/oferta/atrakcyjna-dzialka-budowlana-grajow-wieliczka-ID4vwuO"
